# Update the "Översikt AVESTA" logging-report worksheet:
#  1. Bump the "Förändrad" (Changed) date in column C from 2023-09-13 (45182)
#     to 2023-09-15 (45184) for every existing data row (2..219).
#  2. Restore the explicit 15pt row height on row 219 (lost in a prior edit).
#  3. Append a brand new data row (220) for case "A 43086-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastExistingRow = 219
$newRow = 220

# 1. Update column C ("Förändrad") for all existing data rows to the new date.
$ws.Range("C2:C$lastExistingRow").Value = 45184

# 2. Make sure row 219 carries the explicit custom row height again.
$ws.Rows.Item($lastExistingRow).RowHeight = 15

# 3. Append the new row with the new case data.
$ws.Range("A$newRow").Value = "A 43086-2023"
$ws.Range("B$newRow").Value = 45182
$ws.Range("C$newRow").Value = 45184
$ws.Range("B$newRow`:C$newRow").NumberFormat = "YYYY-MM-DD"
$ws.Range("D$newRow").Value = "DALARNAS LÄN"
$ws.Range("E$newRow").Value = "AVESTA"
$ws.Range("G$newRow").Value = 3.4
$ws.Range("H$newRow`:Q$newRow").Value = 0
$ws.Range("R$newRow").Value = ""
$ws.Range("R$newRow").WrapText = $true
